$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E72 (end time) - this ripples into F72/G72 via the shared formulas
$ws.Range("E72").Value = 0.46180555555555558

# Insert a new row at 74, pushing the old row74:76 (summary rows) down to 75:77,
# and duplicating the blank-placeholder formatting of row 73 into the new row 74.
$ws.Rows("74:74").Insert()

# Turn the former blank placeholder row (73) into a real data row.
$ws.Range("A73").Value = 2014
$ws.Range("B73").Value = 3
$ws.Range("C73").Value = 17
$ws.Range("D73").Value = 0.5625
$ws.Range("E73").Value = 0.625
$ws.Range("F73").Formula = "=(E73-D73)*24*60"
$ws.Range("G73").Formula = "=F73/60"

# Fix up the summary formulas (now on rows 75-77) to include the new blank row 74.
$ws.Range("F75").Formula = "=SUM(F2:F74)"
$ws.Range("F76").Formula = "=F75/60"
$ws.Range("F77").Formula = "=F76/38.5"

# Dimension & view bookkeeping
$ws.Range("F73").Select()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
